# Calibrate Azerbaijan for total population
# - update effective contact rate (B2) and proportion of deaths reported (B3)
# - add two new constant rows: susceptible_fully (starting population) and
#   active (seed of active TB cases), each with a description in column E
# - add a whole-number data validation (0 .. 10,000,000,000) over B4:D5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# --- Update existing parameter values -------------------------------------
$ws.Range("B2").Value = 8.4
$ws.Range("B3").Value = 0.15

# --- Row 4: susceptible_fully ----------------------------------------------
$ws.Range("A4").Value = "susceptible_fully"
$ws.Range("B4").Value = 3240000
$ws.Range("E4").Value = "Starting population of fully susceptible persons (essentially the starting population size)"

# --- Row 5: active -----------------------------------------------------------
$ws.Range("A5").Value = "active"
$ws.Range("B5").Value = 10
$ws.Range("E5").Value = "Seed of patients with active TB"

# --- Data validation: whole numbers between 0 and 10,000,000,000 ----------
$ws.Range("B4:D5").Validation.Add(1, 1, 1, 0, 10000000000)

# --- Selection matches the saved workbook state ----------------------------
$ws.Range("A11").Select()
